$d = $word.ActiveDocument

$pairs = @(
    @("62×91=5642", "19×85=1615"),
    @("98×53=5194", "77×64=4928"),
    @("83×62=5146", "53×91=4823"),
    @("35×73=2555", "62×44=2728"),
    @("42×44=1848", "85×58=4930"),
    @("41×12=492",  "14×40=560"),
    @("26×19=494",  "22×65=1430"),
    @("91×63=5733", "70×53=3710"),
    @("40×96=3840", "73×80=5840"),
    @("21×96=2016", "76×64=4864"),
    @("50×25=1250", "12×96=1152"),
    @("46×36=1656", "36×91=3276"),
    @("83×19=1577", "21×87=1827"),
    @("23×62=1426", "51×34=1734"),
    @("61×87=5307", "16×96=1536"),
    @("89×60=5340", "31×74=2294"),
    @("51×92=4692", "73×80=5840"),
    @("18×99=1782", "89×19=1691"),
    @("77×99=7623", "94×79=7426"),
    @("85×80=6800", "42×51=2142"),
    @("40×76=3040", "90×21=1890"),
    @("44×65=2860", "78×79=6162"),
    @("16×27=432",  "14×26=364"),
    @("65×88=5720", "26×53=1378"),
    @("82×18=1476", "90×24=2160")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
